{"js": "// Office.js (Word JavaScript API) script.\n//\n// The original document was a single paragraph of placeholder filler\n// text (\"One two three ... ten.\" repeated). This rewrites the body as\n// 5 paragraphs of real draft copy for the explanatory text panel:\n//   1) intro paragraph about state borders not matching where people live\n//   2) empty spacer paragraph\n//   3) paragraph introducing the machine-learning re-drawing project\n//   4) empty spacer paragraph\n//   5) '[PANEL OVERVIEW WILL GO HERE]' placeholder paragraph\n//\n// Each target paragraph is expressed below as an array of run strings;\n// they are inserted one after another with insertText(..., \"End\") so the\n// visible text matches exactly (adjoining runs with identical formatting\n// collapse into a single <w:r> on save, same as Word does natively).\nconst PARAS = [\n  [\"States are \", \"a \", \"crucial\", \" \", \"way that the United States divides land into meaningful units.  \", \"Americans deal with different public institutions and are subject to \", \"varying\", \" laws within each state\", \".\", \"  However, state boundaries are perhaps counterintuitive for such a significant set of divisions.  In Kansas City, the Missouri-Kansas boundary cuts the metropolitan area in half, such that Kansas Citians might routinely need to cross state borders to fulfill simple errands.  In Texas, El Paso\", \" is 285 miles from the \", \"nearest\", \" metropolitan area in Texas (Odessa) but only \", \"45\", \" miles from the \", \"closest\", \" \", \"metro\", \" area in New Mexico (Las Cruces).\", \"  In Michigan\\u2019s upper peninsula, Michiganders must cross the Great Lakes over a \", \"five-mile\", \" toll bridge\", \" to reach over 90 percent of Michigan\\u2019s population but share a \", \"200-mile\", \" land border with Wisconsin.\"],\n  [],\n  [\"What if the state borders of the United States matched where Americans live?\", \"  This project applies machine learning \", \"to\", \" \", \"imagine\", \" the United States\\u2019 state borders\", \".  The algorithms \", \"group\", \" people who live near each other and place \", \"boundaries\", \" \", \"in the\", \" large\", \",\", \" unpopulated stretches\", \" of land between \", \"those groups.\"],\n  [],\n  [\"[PANEL OVERVIEW \", \"WILL GO \", \"HERE]\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Seed off the document's existing first paragraph (reused as new paragraph\n// 1 so we don't leave a stray leading empty paragraph behind), and drop any\n// other pre-existing paragraphs -- this source document only ever has the\n// one, but this keeps the script correct even if that ever changes.\nlet current = paragraphs.items[0];\nfor (let i = paragraphs.items.length - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\ncurrent.clear();\nawait context.sync();\n\nfor (let pi = 0; pi < PARAS.length; pi++) {\n  const runs = PARAS[pi];\n\n  if (pi > 0) {\n    // Move on to a freshly inserted paragraph after the one we just filled.\n    current = current.insertParagraph(\"\", \"After\");\n    await context.sync();\n  }\n\n  for (const runText of runs) {\n    current.insertText(runText, \"End\");\n    await context.sync();\n  }\n}\n\nbody.load(\"text\");\nawait context.sync();\nreturn body.text;\n\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Rebuilds the document body as 5 paragraphs:\n#   1) intro paragraph about state borders (built from many runs)\n#   2) empty paragraph\n#   3) paragraph about the machine-learning re-drawing project\n#   4) empty paragraph\n#   5) '[PANEL OVERVIEW WILL GO HERE]' placeholder paragraph\n\n$d = $word.ActiveDocument\n\n# Right single quotation mark (U+2019), used for the possessives below.\n$apos = [char]0x2019\n\n$PARAS = @(\n    @(\"States are \", \"a \", \"crucial\", \" \", \"way that the United States divides land into meaningful units.  \", \"Americans deal with different public institutions and are subject to \", \"varying\", \" laws within each state\", \".\", \"  However, state boundaries are perhaps counterintuitive for such a significant set of divisions.  In Kansas City, the Missouri-Kansas boundary cuts the metropolitan area in half, such that Kansas Citians might routinely need to cross state borders to fulfill simple errands.  In Texas, El Paso\", \" is 285 miles from the \", \"nearest\", \" metropolitan area in Texas (Odessa) but only \", \"45\", \" miles from the \", \"closest\", \" \", \"metro\", \" area in New Mexico (Las Cruces).\", \"  In Michigan\" + $apos + \"s upper peninsula, Michiganders must cross the Great Lakes over a \", \"five-mile\", \" toll bridge\", \" to reach over 90 percent of Michigan\" + $apos + \"s population but share a \", \"200-mile\", \" land border with Wisconsin.\"),\n    @(),\n    @(\"What if the state borders of the United States matched where Americans live?\", \"  This project applies machine learning \", \"to\", \" \", \"imagine\", \" the United States\" + $apos + \" state borders\", \".  The algorithms \", \"group\", \" people who live near each other and place \", \"boundaries\", \" \", \"in the\", \" large\", \",\", \" unpopulated stretches\", \" of land between \", \"those groups.\"),\n    @(),\n    @(\"[PANEL OVERVIEW \", \"WILL GO \", \"HERE]\"),\n)\n\n# Wipe the existing body content (the filler paragraph) before rebuilding.\n$d.Content.Delete()\n\nfunction AppendText($t) {\n    $end = $d.Content.End\n    $r = $d.Range($end - 1, $end - 1)\n    $r.InsertAfter($t)\n}\n\nfunction AppendParaBreak {\n    $end = $d.Content.End\n    $r = $d.Range($end - 1, $end - 1)\n    $r.InsertParagraphAfter()\n}\n\nfor ($i = 0; $i -lt $PARAS.Count; $i++) {\n    if ($i -gt 0) {\n        AppendParaBreak\n    }\n    foreach ($run in $PARAS[$i]) {\n        AppendText $run\n    }\n}\n\n"}
